# Updates the cryptos price-list worksheet: refreshed Price/Volume(1h) figures
# and a few re-ordered coin rows (ranking swaps), per the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.496.92'
$ws.Range("E2").Value = '  +0.97%  '

# Row 3
$ws.Range("D3").Value = '2.432.01'
$ws.Range("E3").Value = '  +0.95%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '''565.92'
$ws.Range("E5").Value = '  +0.71%  '

# Row 6
$ws.Range("D6").Value = '''145.16'
$ws.Range("E6").Value = '  +2.11%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.36%  '

# Row 9
$ws.Range("E9").Value = '  +1.64%  '

# Row 10
$ws.Range("E10").Value = '  +0.43%  '

# Row 11
$ws.Range("D11").Value = '''5.30'
$ws.Range("E11").Value = '  +1.23%  '

# Row 12
$ws.Range("E12").Value = '  +1.90%  '

# Row 13
$ws.Range("D13").Value = '''26.81'
$ws.Range("E13").Value = '  +5.23%  '

# Row 14
$ws.Range("D14").Value = '''0.0000179'
$ws.Range("E14").Value = '  +4.51%  '

# Row 15
$ws.Range("E15").Value = '  +0.62%  '

# Row 16
$ws.Range("D16").Value = '62.409.31'
$ws.Range("E16").Value = '  +0.84%  '

# Row 17
$ws.Range("D17").Value = '2.436.06'
$ws.Range("E17").Value = '  +1.26%  '

# Row 18
$ws.Range("D18").Value = '''11.22'
$ws.Range("E18").Value = '  +0.11%  '

# Row 19
$ws.Range("E19").Value = '  +2.36%  '

# Row 20
$ws.Range("D20").Value = '''323.85'
$ws.Range("E20").Value = '  +0.95%  '

# Row 21
$ws.Range("E21").Value = '  +1.23%  '

# Row 22
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''67.23'
$ws.Range("E23").Value = '  +2.63%  '

# Row 24
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '''1.82'
$ws.Range("E24").Value = '  +5.48%  '

# Row 25
$ws.Range("D25").Value = '''590.48'
$ws.Range("E25").Value = '  +4.68%  '

# Row 26
$ws.Range("D26").Value = '''8.54'
$ws.Range("E26").Value = '  -1.49%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0997'
$ws.Range("E27").Value = '  +7.16%  '

# Row 28
$ws.Range("D28").Value = '2.550.87'
$ws.Range("E28").Value = '  +1.41%  '

# Row 29
$ws.Range("D29").Value = '''8.44'
$ws.Range("E29").Value = '  +3.35%  '

# Row 30
$ws.Range("E30").Value = '  -0.11%  '

# Row 31
$ws.Range("D31").Value = '''1.44'
$ws.Range("E31").Value = '  +4.39%  '

# Row 32
$ws.Range("E32").Value = '  -0.84%  '

# Row 33
$ws.Range("D33").Value = '''1.87'
$ws.Range("E33").Value = '  +0.49%  '

# Row 34
$ws.Range("E34").Value = '  +0.06%  '

# Row 35
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '''4.84'
$ws.Range("E36").Value = '  +1.93%  '

# Row 37
$ws.Range("E37").Value = '  +0.63%  '

# Row 38
$ws.Range("E38").Value = '  +1.48%  '

# Row 39
$ws.Range("D39").Value = '''5.35'
$ws.Range("E39").Value = '  -1.62%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''147.86'
$ws.Range("E40").Value = '  -3.00%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''1.82'
$ws.Range("E41").Value = '  +2.19%  '

# Row 42
$ws.Range("D42").Value = '''1.00'

# Row 43
$ws.Range("D43").Value = '''2.44'
$ws.Range("E43").Value = '  +9.07%  '

# Row 44
$ws.Range("D44").Value = '''148.27'
$ws.Range("E44").Value = '  +0.26%  '

# Row 45
$ws.Range("E45").Value = '  +2.13%  '

# Row 46
$ws.Range("D46").Value = '''0.0535'
$ws.Range("E46").Value = '  +1.50%  '

# Row 47
$ws.Range("D47").Value = '''20.52'
$ws.Range("E47").Value = '  +3.52%  '

# Row 48
$ws.Range("D48").Value = '''0.603'
$ws.Range("E48").Value = '  +2.23%  '

# Row 49
$ws.Range("E49").Value = '  +2.96%  '

# Row 50
$ws.Range("E50").Value = '  +0.22%  '

# Row 51
$ws.Range("E51").Value = '  +4.29%  '
